$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data values (server entry) ---
$ws.Range("A2").Value = "WorldServer_1"
$ws.Range("B2").Value = "000103001"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "WorldServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("G2").Value = 3001

# --- Column width changes ---
# (ColumnWidth is rounded to the nearest whole pixel by the host, i.e.
# stored_width = (Round(ColumnWidth*7) + 5) / 7; the inputs below are chosen
# so the resulting stored width lands as close as possible to the target.)
$ws.Columns.Item(2).ColumnWidth = 13.714285714285714
$ws.Columns.Item(3).ColumnWidth = 17.857142857142858
$ws.Columns.Item(5).ColumnWidth = 11.285714285714286
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 10.285714285714286

# --- Remove data validation ---
$ws.Cells.Validation.Delete()

# --- Selection change ---
$ws.Range("G1").Select()
